$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "2025-06-01 Sunday" "2025-06-02 Monday"

Replace-Text "95×79=" "89×81="
Replace-Text "26×39=" "90×95="
Replace-Text "58×64=" "57×78="
Replace-Text "28×84=" "43×61="
Replace-Text "66×52=" "44×97="
Replace-Text "80×43=" "76×23="
Replace-Text "36×63=" "58×28="
Replace-Text "27×44=" "21×33="
Replace-Text "80×20=" "69×20="
Replace-Text "95×64=" "75×12="
Replace-Text "30×51=" "87×34="
Replace-Text "82×67=" "22×66="
Replace-Text "28×67=" "32×59="
Replace-Text "57×54=" "82×54="
Replace-Text "37×37=" "71×92="
Replace-Text "18×84=" "91×83="
Replace-Text "56×50=" "20×34="
Replace-Text "27×33=" "11×96="
Replace-Text "31×33=" "76×84="
Replace-Text "50×80=" "43×97="
Replace-Text "29×62=" "73×39="
Replace-Text "25×53=" "61×22="
Replace-Text "44×32=" "18×97="
Replace-Text "75×60=" "12×89="
Replace-Text "41×46=" "87×48="
